$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.246.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.32"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5238"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2681"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06366"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.69"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07714"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.602"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.86"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5660"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0" + ([char]8325).ToString() + "8236"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.239.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.709"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.31"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.41"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1203"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.293"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05659"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.277"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.512"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.360"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.583"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.415"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9457"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5791"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.930"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8471"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.026.46"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.42"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.795.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05320"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0" + ([char]8328).ToString() + "103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.045"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.63%  "
